# Appends two new yearly rows (2021, 2022) to Sheet1, matching the existing
# row layout (columns A..AQ), and copies the formatting of column A's
# "year label" cell (bold, centered, thin-bordered) onto the two new labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2021 = @("2021年", 3122.96, 1022.22, 174.75, 0.44, 2240.45, 8222.42, 683.08, 6430.68, 493.28, 460.25, 408.11, -45.3, 852.4400000000001, 3591.43, 536.28, 482.51, 1838.48, 452.39, 5646.29, 1188.09, 7168.33, 908.88, 1935.95, 4756.5, 674.74, 2738.32, 1627.02, 1346.34, 883.92, 92933.03, 9018.57, 3356.11, 958.12, 2771.18, 73.62, 2532.38, 714.95, 6032.2, 495.73, 1738.87, 4567.01, 833.04)
$row2022 = @("2022年", 3056.4, 1017.6, 112.2, 0.6, 1901.1, 7302.6, 241.3, 4288.7, 431.3, 471.2, 307.5, 6.7, 772.2, 2571.5, 743.5, 436.6, 1644, 397.6, 5319.6, 1323.1, 10202, 763.1, 3154, 5915.6, 614.4, 451.4, 3545, 1000.8, 763.8, 84038.5, 7389.5, 3250.3, 621.1, 3116.3, 73, 2065.7, 768.7, 4759, 480.9, 1797.9, 365.5, 594.9)

$newRows = @(11, 12)
$rowData = @($row2021, $row2022)

for ($r = 0; $r -lt $newRows.Length; $r++) {
    $targetRow = $newRows[$r]
    $values = $rowData[$r]
    for ($c = 1; $c -le $values.Length; $c++) {
        $ws.Cells.Item($targetRow, $c).Value = $values[$c - 1]
    }
}

# Copy the formatting (bold, centered, top-aligned, thin border) of the
# year-label column from the last existing row (A10) onto the two new
# year-label cells (A11, A12), matching the style used by A2:A10.
$ws.Range("A10").Copy()
$ws.Range("A11:A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0
